# SESUserManage.xlsx — approve the Test-Cases row 5 entry and move the
# active selection to I5 (scrolling the view left so column B is the
# left-most visible column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# --- content change --------------------------------------------------
# I5 ("Approved/Rejected") flips from "Rejected" to "Approved".
$ws.Range("I5").Value = "Approved"

# --- view/selection change --------------------------------------------
# Bring the sheet to the front, scroll so column B is left-most
# (was F1), and move the single-cell selection to I5 (was I2:I3).
$ws.Activate()
$ws.Range("I5").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
